# Week 6 update: insert a new "parking" slide right before the
# "T-tests are powerful" slide (old slide #20) to mark where lecture
# stopped, per commit message "sending week 6 updates - distributions".

$p = $ppt.ActivePresentation

# "Title and Content" custom layout (same layout used by the
# surrounding T-test / Shapiro-Wilk slides).
$layout = $p.SlideMaster.CustomLayouts.Item(2)

# Insert the new slide at position 20, pushing "T-tests are powerful"
# (and everything after it) down by one.
$newSlide = $p.Slides.AddSlide(20, $layout)

# Leave the title placeholder empty and set the body placeholder text.
$newSlide.Shapes.Item(2).TextFrame.TextRange.Text = "STOPPED HERE!!"
